$wb = $excel.ActiveWorkbook

# --- Update "Logs" sheet: append new row with the latest e-mail ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A9").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B9").Value = "mailmind.test@zohomail.eu"
$logs.Range("C9").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D9").Value = "Offerte / Prijsaanvraag"
$logs.Range("F9").Value = "2025-06-19 21:18:10"
$logs.Range("G9").Value = "Nee"

# --- Update "Dashboard" sheet: append the new category tally ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A8").Value = "Offerte / Prijsaanvraag"
$dash.Range("B8").Value = 1

# --- Extend the chart's category/value source ranges to include row 8 ---
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$8,'Dashboard'!`$B`$2:`$B`$8,1)"

# --- Grow the conditional-formatting ranges on "Logs" so they cover the new row 9 ---
$catFmt = $logs.Range("D2:D8").FormatConditions
for ($i = 1; $i -le $catFmt.Count; $i++) {
    $catFmt.Item($i).ModifyAppliesToRange($logs.Range("D2:D9"))
}

$answeredFmt = $logs.Range("G2:G8").FormatConditions
for ($i = 1; $i -le $answeredFmt.Count; $i++) {
    $answeredFmt.Item($i).ModifyAppliesToRange($logs.Range("G2:G9"))
}
